# Auto-generated edit script applying the scheduled-runner update to Sheets/Omega_Profits.xlsx
# Updates currentAveragePrice / LevePrice / LeveProfit figures across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 18114.234
$ws.Range("J51").Value = 18629.467
$ws.Range("L51").Value = 18629.467
$ws.Range("N51").Value = -19597.467
$ws.Range("H70").Value = 12830.167
$ws.Range("I70").Value = 991
$ws.Range("K70").Value = 2973
$ws.Range("M70").Value = -2703
$ws.Range("H73").Value = 12830.167
$ws.Range("I73").Value = 991
$ws.Range("K73").Value = 2973
$ws.Range("M73").Value = -2037
$ws.Range("H107").Value = 1406.875
$ws.Range("I107").Value = 1005.0833
$ws.Range("K107").Value = 1005.0833
$ws.Range("M107").Value = 914.9167
$ws.Range("H125").Value = 1540
$ws.Range("I125").Value = 1471
$ws.Range("J125").Value = 1712.5
$ws.Range("K125").Value = 13239
$ws.Range("L125").Value = 15412.5
$ws.Range("M125").Value = -10779
$ws.Range("N125").Value = -20332.5

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 14987.412
$ws.Range("I45").Value = 18406.924
$ws.Range("K45").Value = 18406.924
$ws.Range("M45").Value = -18029.924
$ws.Range("H61").Value = 7833.5713
$ws.Range("I61").Value = 6960.25
$ws.Range("K61").Value = 6960.25
$ws.Range("M61").Value = -6748.25
$ws.Range("H110").Value = 3059.6667
$ws.Range("I110").Value = 2880
$ws.Range("J110").Value = 3149.5
$ws.Range("K110").Value = 2880
$ws.Range("L110").Value = 3149.5
$ws.Range("M110").Value = -835
$ws.Range("N110").Value = -7239.5
$ws.Range("H122").Value = 2425.8
$ws.Range("I122").Value = 2452.8462
$ws.Range("K122").Value = 7358.5386
$ws.Range("M122").Value = -4908.5386
$ws.Range("H132").Value = 3283
$ws.Range("I132").Value = 2855.125
$ws.Range("K132").Value = 8565.375
$ws.Range("M132").Value = -6035.375
$ws.Range("H134").Value = 64998.75
$ws.Range("J134").Value = 64998.75
$ws.Range("L134").Value = 64998.75
$ws.Range("N134").Value = -75138.75
$ws.Range("H136").Value = 7833.5713
$ws.Range("I136").Value = 6960.25
$ws.Range("K136").Value = 20880.75
$ws.Range("M136").Value = -18330.75

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 77589.8
$ws.Range("J132").Value = 77589.8
$ws.Range("L132").Value = 77589.8
$ws.Range("N132").Value = -87709.8

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 18499.75
$ws.Range("I39").Value = 11333.333
$ws.Range("K39").Value = 11333.333
$ws.Range("M39").Value = -10942.333
$ws.Range("H49").Value = 18499.75
$ws.Range("I49").Value = 11333.333
$ws.Range("K49").Value = 11333.333
$ws.Range("M49").Value = -11151.333
$ws.Range("H86").Value = 9480
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()  # was -98877
$ws.Range("H89").Value = 9480
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()  # was -494384
$ws.Range("H99").Value = 5398.2856
$ws.Range("J99").Value = 7138.4
$ws.Range("L99").Value = 7138.4
$ws.Range("N99").Value = -10134.4
$ws.Range("H126").Value = 5398.2856
$ws.Range("J126").Value = 7138.4
$ws.Range("L126").Value = 21415.2
$ws.Range("N126").Value = -26355.2
$ws.Range("H135").Value = 79811.78
$ws.Range("J135").Value = 80949.625
$ws.Range("L135").Value = 80949.625
$ws.Range("N135").Value = -91089.625

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1749.1111
$ws.Range("J11").Value = 3574.75
$ws.Range("L11").Value = 10724.25
$ws.Range("N11").Value = -11004.25
$ws.Range("H33").Value = 146
$ws.Range("I33").Value = 65.5
$ws.Range("J33").Value = 199.66667
$ws.Range("K33").Value = 393
$ws.Range("L33").Value = 1198.00002
$ws.Range("M33").Value = -110
$ws.Range("N33").Value = -1764.00002
$ws.Range("H75").Value = 3849.4
$ws.Range("I75").Value = 3846
$ws.Range("J75").Value = 3850.25
$ws.Range("K75").Value = 11538
$ws.Range("L75").Value = 11550.75
$ws.Range("M75").Value = -10540
$ws.Range("N75").Value = -13546.75
$ws.Range("H78").Value = 3849.4
$ws.Range("I78").Value = 3846
$ws.Range("J78").Value = 3850.25
$ws.Range("K78").Value = 34614
$ws.Range("L78").Value = 34652.25
$ws.Range("M78").Value = -29622
$ws.Range("N78").Value = -44636.25
$ws.Range("H97").Value = 999
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()  # was -5492
$ws.Range("H117").Value = 2865.7273
$ws.Range("J117").Value = 3215.375
$ws.Range("L117").Value = 9646.125
$ws.Range("N117").Value = -16530.125
$ws.Range("H120").Value = 10202
$ws.Range("I120").Value = 3603
$ws.Range("J120").Value = 29999
$ws.Range("K120").Value = 10809
$ws.Range("L120").Value = 89997
$ws.Range("M120").Value = -5971
$ws.Range("N120").Value = -99673
$ws.Range("H140").Value = 3587.3
$ws.Range("I140").Value = 2609.6875
$ws.Range("K140").Value = 7829.0625
$ws.Range("M140").Value = -2649.0625
$ws.Range("H141").Value = 4101.625
$ws.Range("I141").Value = 4101.625
$ws.Range("K141").Value = 12304.875
$ws.Range("M141").Value = -7124.875

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3765410.2
$ws.Range("I11").Value = 150675.7
$ws.Range("K11").Value = 150675.7
$ws.Range("M11").Value = -150536.7
$ws.Range("H80").Value = 4960.909
$ws.Range("I80").Value = 3699.25
$ws.Range("J80").Value = 8325.333000000001
$ws.Range("K80").Value = 3699.25
$ws.Range("L80").Value = 8325.333000000001
$ws.Range("M80").Value = -2701.25
$ws.Range("N80").Value = -10321.333
$ws.Range("H83").Value = 4960.909
$ws.Range("I83").Value = 3699.25
$ws.Range("J83").Value = 8325.333000000001
$ws.Range("K83").Value = 18496.25
$ws.Range("L83").Value = 41626.665
$ws.Range("M83").Value = -13504.25
$ws.Range("N83").Value = -51610.665
$ws.Range("H102").Value = 2613.2856
$ws.Range("I102").Value = 2613.2856
$ws.Range("K102").Value = 2613.2856
$ws.Range("M102").Value = -991.2856000000002
$ws.Range("H107").Value = 510.69232
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()  # was -5040
$ws.Range("H126").Value = 8642.571
$ws.Range("I126").Value = 7999.5
$ws.Range("J126").Value = 8899.799999999999
$ws.Range("K126").Value = 23998.5
$ws.Range("L126").Value = 26699.4
$ws.Range("M126").Value = -21528.5
$ws.Range("N126").Value = -31639.4
$ws.Range("H132").Value = 5415.2334
$ws.Range("I132").Value = 5617.72
$ws.Range("K132").Value = 16853.16
$ws.Range("M132").Value = -14323.16

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11187.704
$ws.Range("I7").Value = 11251.087
$ws.Range("K7").Value = 11251.087
$ws.Range("M7").Value = -11139.087
$ws.Range("H16").Value = 819.4516
$ws.Range("I16").Value = 665.65515
$ws.Range("K16").Value = 665.65515
$ws.Range("M16").Value = -495.65515
$ws.Range("H40").Value = 4298
$ws.Range("I40").Value = 1384.6666
$ws.Range("K40").Value = 1384.6666
$ws.Range("M40").Value = -1248.6666
$ws.Range("H120").Value = 81999.5
$ws.Range("J120").Value = 81999.5
$ws.Range("L120").Value = 81999.5
$ws.Range("N120").Value = -91675.5
$ws.Range("H126").Value = 11187.704
$ws.Range("I126").Value = 11251.087
$ws.Range("K126").Value = 33753.261
$ws.Range("M126").Value = -31283.261
$ws.Range("H132").Value = 3132.7273
$ws.Range("I132").Value = 3444
$ws.Range("J132").Value = 3063.5557
$ws.Range("K132").Value = 10332
$ws.Range("L132").Value = 9190.667099999999
$ws.Range("M132").Value = -7802
$ws.Range("N132").Value = -14250.6671
$ws.Range("H133").Value = 82499
$ws.Range("J133").Value = 82499
$ws.Range("L133").Value = 82499
$ws.Range("N133").Value = -87559
$ws.Range("H135").Value = 86874
$ws.Range("J135").Value = 86874
$ws.Range("L135").Value = 86874
$ws.Range("N135").Value = -97014

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7820.3335
$ws.Range("I41").Value = 7820.3335
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 7820.3335
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -7430.3335
$ws.Range("N41").ClearContents()  # was -10013
$ws.Range("H107").Value = 3998.75
$ws.Range("J107").Value = 7633.25
$ws.Range("L107").Value = 22899.75
$ws.Range("N107").Value = -26739.75
$ws.Range("H122").Value = 3212.375
$ws.Range("I122").Value = 2966.5
$ws.Range("J122").Value = 3950
$ws.Range("K122").Value = 8899.5
$ws.Range("L122").Value = 11850
$ws.Range("M122").Value = -6449.5
$ws.Range("N122").Value = -16750
$ws.Range("H126").Value = 2290.3635
$ws.Range("I126").Value = 2266.6667
$ws.Range("J126").Value = 2397
$ws.Range("K126").Value = 6800.000100000001
$ws.Range("L126").Value = 7191
$ws.Range("M126").Value = -4330.000100000001
$ws.Range("N126").Value = -12131
$ws.Range("H132").Value = 4369.364
$ws.Range("I132").Value = 4518.1665
$ws.Range("K132").Value = 13554.4995
$ws.Range("M132").Value = -11024.4995
